# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition): update "want to go" counts ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 745
$ws1.Range("F3").Value = 4131
$ws1.Range("F4").Value = 118
$ws1.Range("F5").Value = 751

# --- Sheet "演出" (Performance): add a new event row ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("A2").Value = 1
$ws2.Range("A2").Style = $ws2.Range("A1").Style
$ws2.Range("B2").Value = "2024-06-22"
$ws2.Range("C2").Value = "南宁·浪漫古典·百年经典世界名曲音乐会"
$ws2.Range("D2").Value = "广西壮族自治区南宁市良庆区龙堤路25号  广西文化艺术中心-音乐厅"
$ws2.Range("E2").Value = "2024.06.22 20:00-06.22 21:30"
$ws2.Range("F2").Value = 0
$ws2.Range("G2").Value = 50
$ws2.Range("H2").Value = "https://show.bilibili.com/platform/detail.html?id=83959"
$ws2.Range("I2").Value = "//i1.hdslb.com/bfs/openplatform/202404/H0f8U7no1712041461015.jpeg"

# --- Sheet "全部类型" (All types): same count updates, plus new event row ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 745
$ws4.Range("F3").Value = 4131
$ws4.Range("F4").Value = 118
$ws4.Range("F5").Value = 751

$ws4.Range("A6").Value = 5
$ws4.Range("A6").Style = $ws4.Range("A1").Style
$ws4.Range("B6").Value = "2024-06-22"
$ws4.Range("C6").Value = "南宁·浪漫古典·百年经典世界名曲音乐会"
$ws4.Range("D6").Value = "广西壮族自治区南宁市良庆区龙堤路25号  广西文化艺术中心-音乐厅"
$ws4.Range("E6").Value = "2024.06.22 20:00-06.22 21:30"
$ws4.Range("F6").Value = 0
$ws4.Range("G6").Value = 50
$ws4.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=83959"
$ws4.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202404/H0f8U7no1712041461015.jpeg"
